$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(5, 6).Value = 67
$ws.Cells.Item(6, 6).Value = 825
$ws.Cells.Item(7, 6).Value = 404
$ws.Cells.Item(8, 6).Value = 4674
$ws.Cells.Item(9, 6).Value = 4674
$ws.Cells.Item(11, 6).Value = 120
$ws.Cells.Item(12, 6).Value = 154
$ws.Cells.Item(15, 6).Value = 113
$ws.Cells.Item(16, 6).Value = 7419
$ws.Cells.Item(17, 6).Value = 247
$ws.Cells.Item(18, 6).Value = 125
$ws.Cells.Item(20, 6).Value = 22
$ws.Cells.Item(21, 6).Value = 515
$ws.Cells.Item(22, 6).Value = 1346
$ws.Cells.Item(25, 6).Value = 1738
$ws.Cells.Item(27, 6).Value = 1990
$ws.Cells.Item(28, 6).Value = 6158
$ws.Cells.Item(30, 6).Value = 21
$ws.Cells.Item(31, 6).Value = 114
$ws.Cells.Item(33, 6).Value = 444
$ws.Cells.Item(34, 6).Value = 6381
$ws.Cells.Item(35, 6).Value = 23
$ws.Cells.Item(36, 6).Value = 205
$ws.Cells.Item(39, 6).Value = 19
$ws.Cells.Item(41, 6).Value = 2452
$ws.Cells.Item(43, 6).Value = 56
$ws.Cells.Item(44, 6).Value = 1016
$ws.Cells.Item(45, 6).Value = 37
$ws.Cells.Item(46, 6).Value = 426
$ws.Cells.Item(47, 6).Value = 2133
$ws.Cells.Item(48, 6).Value = 41
$ws.Cells.Item(49, 6).Value = 1071

$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(3, 6).Value = 229
$ws.Cells.Item(4, 6).Value = 8
$ws.Cells.Item(6, 6).Value = 120
$ws.Cells.Item(8, 6).Value = 10
$ws.Cells.Item(10, 6).Value = 7

$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(2, 6).Value = 1441

$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(3, 6).Value = 1441
$ws.Cells.Item(5, 6).Value = 229
$ws.Cells.Item(6, 6).Value = 67
$ws.Cells.Item(8, 6).Value = 404
$ws.Cells.Item(9, 6).Value = 4674
$ws.Cells.Item(10, 6).Value = 4674
$ws.Cells.Item(12, 6).Value = 120
$ws.Cells.Item(13, 6).Value = 154
$ws.Cells.Item(16, 6).Value = 113
$ws.Cells.Item(17, 6).Value = 7419
$ws.Cells.Item(18, 6).Value = 247
$ws.Cells.Item(19, 6).Value = 125
$ws.Cells.Item(20, 6).Value = 515
$ws.Cells.Item(21, 6).Value = 1346
$ws.Cells.Item(22, 6).Value = 120
$ws.Cells.Item(24, 6).Value = 1738
$ws.Cells.Item(26, 6).Value = 1990
$ws.Cells.Item(27, 6).Value = 10
$ws.Cells.Item(29, 6).Value = 6158
$ws.Cells.Item(31, 6).Value = 7
$ws.Cells.Item(32, 6).Value = 21
$ws.Cells.Item(33, 6).Value = 114
$ws.Cells.Item(35, 6).Value = 444
$ws.Cells.Item(36, 6).Value = 6381
$ws.Cells.Item(37, 6).Value = 23
$ws.Cells.Item(38, 6).Value = 205
$ws.Cells.Item(40, 6).Value = 19
$ws.Cells.Item(42, 6).Value = 2452
$ws.Cells.Item(44, 6).Value = 1016
$ws.Cells.Item(45, 6).Value = 37
$ws.Cells.Item(46, 6).Value = 426
$ws.Cells.Item(48, 6).Value = 2133
$ws.Cells.Item(49, 6).Value = 41

